# Update "想去人数" (F column) figures for a handful of events.
# These edits appear on both the "展览" sheet and the aggregated
# "全部类型" sheet, which duplicates the same rows.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row -> new F value
$exhibition = $wb.Worksheets.Item("展览")
$exhibition.Range("F2").Value  = 4561
$exhibition.Range("F3").Value  = 2502
$exhibition.Range("F10").Value = 174
$exhibition.Range("F11").Value = 174
$exhibition.Range("F12").Value = 1705
$exhibition.Range("F13").Value = 311
$exhibition.Range("F14").Value = 3755
$exhibition.Range("F15").Value = 25
$exhibition.Range("F16").Value = 250

# Sheet "全部类型" (all types, combines entries from every category)
$allTypes = $wb.Worksheets.Item("全部类型")
$allTypes.Range("F2").Value  = 4561
$allTypes.Range("F3").Value  = 2502
$allTypes.Range("F12").Value = 174
$allTypes.Range("F13").Value = 174
$allTypes.Range("F16").Value = 1705
$allTypes.Range("F17").Value = 311
$allTypes.Range("F18").Value = 3755
$allTypes.Range("F19").Value = 25
$allTypes.Range("F20").Value = 250
